$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 rows to make room for the new "M2" sending-cluster block
# (originally rows 11:13 held the "sCs" block; push it down to 14:16)
$ws.Range("A11:T13").EntireRow.Insert()

$rng = $ws.Range("A2:T2")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 55.26722666666667
$arr[0,7] = 165.80168
$arr[0,8] = 0.01597985502890189
$arr[0,9] = 0.01623411501809385
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 5.495057333333333
$arr[0,13] = 16.485172
$arr[0,14] = 0.8161989011161211
$arr[0,15] = 0.8403205285996808
$arr[0,16] = 303.6965791876622
$arr[0,17] = 2733.26921268896
$arr[0,18] = 0.01304274011458465
$arr[0,19] = 0.01364186011335264
$rng.Value = $arr

$rng = $ws.Range("A3:T3")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 55.26722666666667
$arr[0,7] = 165.80168
$arr[0,8] = 0.01597985502890189
$arr[0,9] = 0.01623411501809385
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.657666
$arr[0,13] = 1.972998
$arr[0,14] = 0.09768528951377062
$arr[0,15] = 0.1005722428790014
$arr[0,16] = 36.34737589296
$arr[0,17] = 327.12638303664
$arr[0,18] = 0.001560996764886365
$arr[0,19] = 0.001632701358525379
$rng.Value = $arr

$rng = $ws.Range("A4:T4")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "sCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 55.26722666666667
$arr[0,7] = 165.80168
$arr[0,8] = 0.01597985502890189
$arr[0,9] = 0.01623411501809385
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 0.5797745
$arr[0,13] = 1.159549
$arr[0,14] = 0.08611580937010824
$arr[0,15] = 0.0591072285213179
$arr[0,16] = 32.04252870705333
$arr[0,17] = 192.25517224232
$arr[0,18] = 0.001376118149430881
$arr[0,19] = 0.0009595535462158321
$rng.Value = $arr

$rng = $ws.Range("A5:T5")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 3231.170084666667
$arr[0,7] = 9693.510254
$arr[0,8] = 0.93425403518284
$arr[0,9] = 0.9491192151521513
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 5.495057333333333
$arr[0,13] = 16.485172
$arr[0,14] = 0.8161989011161211
$arr[0,15] = 0.8403205285996808
$arr[0,16] = 17755.46486899485
$arr[0,17] = 159799.1838209537
$arr[0,18] = 0.762537116879536
$arr[0,19] = 0.7975643605807698
$rng.Value = $arr

$rng = $ws.Range("A6:T6")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 3231.170084666667
$arr[0,7] = 9693.510254
$arr[0,8] = 0.93425403518284
$arr[0,9] = 0.9491192151521513
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.657666
$arr[0,13] = 1.972998
$arr[0,14] = 0.09768528951377062
$arr[0,15] = 0.1005722428790014
$arr[0,16] = 2125.030704902388
$arr[0,17] = 19125.27634412149
$arr[0,18] = 0.09126287590624417
$arr[0,19] = 0.09545504822740933
$rng.Value = $arr

$rng = $ws.Range("A7:T7")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "sCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 3231.170084666667
$arr[0,7] = 9693.510254
$arr[0,8] = 0.93425403518284
$arr[0,9] = 0.9491192151521513
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 0.5797745
$arr[0,13] = 1.159549
$arr[0,14] = 0.08611580937010824
$arr[0,15] = 0.0591072285213179
$arr[0,16] = 1873.350020252574
$arr[0,17] = 11240.10012151545
$arr[0,18] = 0.08045404239705985
$arr[0,19] = 0.0560998063439721
$rng.Value = $arr

$rng = $ws.Range("A8:T8")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "M1"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.924789000000001
$arr[0,7] = 17.774367
$arr[0,8] = 0.001713081603820286
$arr[0,9] = 0.001740338929326963
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 5.495057333333333
$arr[0,13] = 16.485172
$arr[0,14] = 0.8161989011161211
$arr[0,15] = 0.8403205285996808
$arr[0,16] = 32.55705524290266
$arr[0,17] = 293.013497186124
$arr[0,18] = 0.00139821532256036
$arr[0,19] = 0.001462442529034636
$rng.Value = $arr

$rng = $ws.Range("A9:T9")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "M1"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "FAPs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.924789000000001
$arr[0,7] = 17.774367
$arr[0,8] = 0.001713081603820286
$arr[0,9] = 0.001740338929326963
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.657666
$arr[0,13] = 1.972998
$arr[0,14] = 0.09768528951377062
$arr[0,15] = 0.1005722428790014
$arr[0,16] = 3.896532282474
$arr[0,17] = 35.068790542266
$arr[0,18] = 0.0001673428724298992
$arr[0,19] = 0.0001750297894920526
$rng.Value = $arr

$rng = $ws.Range("A10:T10")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "M1"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "sCs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 5.924789000000001
$arr[0,7] = 17.774367
$arr[0,8] = 0.001713081603820286
$arr[0,9] = 0.001740338929326963
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 0.5797745
$arr[0,13] = 1.159549
$arr[0,14] = 0.08611580937010824
$arr[0,15] = 0.0591072285213179
$arr[0,16] = 3.4350415800805
$arr[0,17] = 20.610249480483
$arr[0,18] = 0.0001475234088300271
$arr[0,19] = 0.0001028666108002745
$rng.Value = $arr

$rng = $ws.Range("A11:T11")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "M2"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 3.689559
$arr[0,7] = 11.068677
$arr[0,8] = 0.001066791686439732
$arr[0,9] = 0.001083765710432669
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 5.495057333333333
$arr[0,13] = 16.485172
$arr[0,14] = 0.8161989011161211
$arr[0,15] = 0.8403205285996808
$arr[0,16] = 20.274338239716
$arr[0,17] = 182.469044157444
$arr[0,18] = 0.0008707142021919227
$arr[0,19] = 0.0009107105746689886
$rng.Value = $arr

$rng = $ws.Range("A12:T12")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "M2"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "FAPs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 3.689559
$arr[0,7] = 11.068677
$arr[0,8] = 0.001066791686439732
$arr[0,9] = 0.001083765710432669
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.657666
$arr[0,13] = 1.972998
$arr[0,14] = 0.09768528951377062
$arr[0,15] = 0.1005722428790014
$arr[0,16] = 2.426497509294
$arr[0,17] = 21.838477583646
$arr[0,18] = 0.0001042098547407488
$arr[0,19] = 0.0001089967482535678
$rng.Value = $arr

$rng = $ws.Range("A13:T13")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "M2"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "sCs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 3.689559
$arr[0,7] = 11.068677
$arr[0,8] = 0.001066791686439732
$arr[0,9] = 0.001083765710432669
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 0.5797745
$arr[0,13] = 1.159549
$arr[0,14] = 0.08611580937010824
$arr[0,15] = 0.0591072285213179
$arr[0,16] = 2.1391122244455
$arr[0,17] = 12.834673346673
$arr[0,18] = 0.00009186762950706021
$arr[0,19] = 0.00006405838751011219
$rng.Value = $arr

$rng = $ws.Range("A14:T14")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "sCs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 1
$arr[0,6] = 162.5045395
$arr[0,7] = 325.009079
$arr[0,8] = 0.046986236497998
$arr[0,9] = 0.03182256518999536
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 5.495057333333333
$arr[0,13] = 16.485172
$arr[0,14] = 0.8161989011161211
$arr[0,15] = 0.8403205285996808
$arr[0,16] = 892.9717614794313
$arr[0,17] = 5357.830568876589
$arr[0,18] = 0.03835011459724815
$arr[0,19] = 0.0267411548018547
$rng.Value = $arr

$rng = $ws.Range("A15:T15")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "sCs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "FAPs"
$arr[0,4] = 2
$arr[0,5] = 1
$arr[0,6] = 162.5045395
$arr[0,7] = 325.009079
$arr[0,8] = 0.046986236497998
$arr[0,9] = 0.03182256518999536
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 0.657666
$arr[0,13] = 1.972998
$arr[0,14] = 0.09768528951377062
$arr[0,15] = 0.1005722428790014
$arr[0,16] = 106.873710474807
$arr[0,17] = 641.2422628488421
$arr[0,18] = 0.00458986411546943
$arr[0,19] = 0.003200466755321069
$rng.Value = $arr

$rng = $ws.Range("A16:T16")
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "sCs"
$arr[0,1] = "Bgn"
$arr[0,2] = "Fgfr3"
$arr[0,3] = "sCs"
$arr[0,4] = 2
$arr[0,5] = 1
$arr[0,6] = 162.5045395
$arr[0,7] = 325.009079
$arr[0,8] = 0.046986236497998
$arr[0,9] = 0.03182256518999536
$arr[0,10] = 2
$arr[0,11] = 1
$arr[0,12] = 0.5797745
$arr[0,13] = 1.159549
$arr[0,14] = 0.08611580937010824
$arr[0,15] = 0.0591072285213179
$arr[0,16] = 94.21598813634276
$arr[0,17] = 376.863952545371
$arr[0,18] = 0.004046257785280417
$arr[0,19] = 0.001880943632819592
$rng.Value = $arr

